$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell updates derived from the crypto price refresh diff.
# Price cells (column D) contain digit-and-dot strings that Excel would
# otherwise auto-convert to numbers, so we force text format, assign, then
# restore the default "Normal" style to avoid leaving stray formatting.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.181.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.807.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "339.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3936"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3487"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.20"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.167"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07543"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.805.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.142"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001100"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06708"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("E21").Value = "  +1.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.554"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.152.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.398"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.32%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.478"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.511"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "154.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.011.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "135.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.184"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.028"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08845"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02427"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.45%  "
$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6907"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06523"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.435"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.610"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2206"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.259"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.483"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.873"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.141"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "130.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07180"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.32%  "
